$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "global_context.js" (2nd sheet): replace the first three
# key/value rows (global.key1/2/3) with mail settings used by the new
# HtmReporter / EmailAlerter (mail_from, mail_to, mail_subject).
# ------------------------------------------------------------------
$wsGlobal = $wb.Worksheets.Item(2)
$wsInit   = $wb.Worksheets.Item(3)

# Drop the old global.key1 / global.key2 / global.key3 rows (rows 2-4),
# then make room for the three new mail rows.
$wsGlobal.Rows("2:4").Delete()
$wsGlobal.Rows("2:2").Insert()
$wsGlobal.Rows("2:2").Insert()
$wsGlobal.Rows("2:2").Insert()

# global.mail_from (row 3 first, to reproduce original shared-string order)
$wsGlobal.Range("A3").Value = "global.mail_from"

# "init.js" sheet gains a line that stashes global_context for later use.
$wsInit.Range("A3").Value = "env_context[""global_context""] = global_context;"

# global.mail_to / global.mail_subject + the two email values
$wsGlobal.Range("A2").Value = "global.mail_to"
$wsGlobal.Range("B2").Value = "2476382757@qq.com"
$wsGlobal.Range("A4").Value = "global.mail_subject"
$wsGlobal.Range("B4").Value = "自动测试报告"
$wsGlobal.Range("B3").Value = "2476382757@my-home-comp.com"

# Hyperlink the two mail addresses (mailto:) - mirrors how Excel turns
# an email address typed into a cell into a clickable hyperlink.
$wsGlobal.Hyperlinks.Add($wsGlobal.Range("B2"), "mailto:2476382757@qq.com")
$wsGlobal.Hyperlinks.Add($wsGlobal.Range("B3"), "mailto:2476382757@my-home-comp.com")

# Column widths grew to fit the new, longer content.
$wsGlobal.Columns("A").ColumnWidth = 20.857142857142858
$wsGlobal.Columns("B").ColumnWidth = 18.714285714285715

$wsInit.Columns("A").ColumnWidth = 52

# Selection left on init.js at the newly-added cell.
$wsInit.Range("A3").Select()

# Finally activate the global_context.js tab and leave the selection on
# B4 (global.mail_subject's value), matching the saved workbook state.
$wsGlobal.Activate()
$wsGlobal.Range("B4").Select()
